$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells ("BL" and "Operating Freq") in columns F and G
$ws.Cells.Item(1, 6).Value = "BL"
$ws.Cells.Item(1, 7).Value = "Operating Freq"

# Give the new header cells the same (centered) formatting as the other
# header cells in row 1 by copying the format from A1.
$ws.Range("A1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

# Add the new numeric data cells for rows 2 and 3 in columns F and G
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 0

# Move the active selection to G4, matching the saved view state
$ws.Range("G4").Select() | Out-Null
